$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.7357868552207947
$ws.Range("B1").Value = 2.255160331726074
$ws.Range("C1").Value = 3.868005752563477
$ws.Range("D1").Value = 3.42602014541626
$ws.Range("E1").Value = 1.999013423919678
